$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Relabel polymer names in Sheet1: strip the "DIP " prefix, and further
# shorten "S1" -> "S" and "B1" -> "B" for the single-sample categories.
$ws1.Range("A3").Value = "S"
$ws1.Range("A4").Value = "B"
$ws1.Range("A5").Value = "G1"
$ws1.Range("A6").Value = "G2"
$ws1.Range("A7").Value = "G3"

# Update the active selection to reflect where the edit left the cursor.
$ws1.Range("A8").Select()
